# Insert a new data row at row 63 (pushes the existing rows 63..180 down
# to 64..181, growing the used range from A1:R180 to A1:R181) and fill the
# newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63..180 down by one, carrying formatting (date style on col D)
# from the row above, exactly like Excel's native Insert does.
$ws.Rows.Item(63).Insert()

# Populate the freshly inserted row 63 with the new record's data.
$ws.Cells.Item(63, 1).Value2 = 10
$ws.Cells.Item(63, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value2 = "La Araucanía"
$ws.Cells.Item(63, 4).Value2 = 44469
$ws.Cells.Item(63, 5).Value2 = 9
$ws.Cells.Item(63, 6).Value2 = 100112017
$ws.Cells.Item(63, 7).Value2 = "Apio"
$ws.Cells.Item(63, 8).Value2 = "Americana (o)"
$ws.Cells.Item(63, 9).Value2 = "Primera"
$ws.Cells.Item(63, 10).Value2 = 240
$ws.Cells.Item(63, 11).Value2 = 9000
$ws.Cells.Item(63, 12).Value2 = 10000
$ws.Cells.Item(63, 13).Value2 = 9417
$ws.Cells.Item(63, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(63, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(63, 16).Value2 = 1570
$ws.Cells.Item(63, 17).Value2 = 6
$ws.Cells.Item(63, 18).Value2 = "Hortaliza"
